# Update the "想去人数" (interested-count) figures in column F on the
# "展览" sheet (1st sheet) and the "全部类型" sheet (4th sheet), matching
# the refreshed data pulled from bilibili at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item(1)   # 展览
$sheetAllTypes   = $wb.Worksheets.Item(4)   # 全部类型

# Row -> new F value for the "展览" sheet
$exhibitionUpdates = @{
    5  = 913
    7  = 922
    8  = 718
    9  = 169
    11 = 78
    12 = 759
    13 = 249
    14 = 541
    15 = 484
    16 = 1277
    18 = 416
    19 = 1042
    20 = 2772
    21 = 1253
    22 = 639
    24 = 1234
    26 = 958
    27 = 314
    28 = 477
    29 = 1295
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row -> new F value for the "全部类型" sheet
$allTypesUpdates = @{
    12 = 913
    15 = 922
    16 = 718
    17 = 170
    23 = 78
    25 = 759
    26 = 249
    27 = 541
    28 = 484
    29 = 1277
    31 = 416
    32 = 1042
    33 = 2772
    34 = 1253
    35 = 639
    37 = 1234
    40 = 958
    41 = 314
    42 = 477
    43 = 1295
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
